$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while preserving it as TEXT (the workbook
# stores these coin prices / percentages as inline strings, not numbers) and
# without leaving a lingering Text number-format on the cell. We briefly force
# the cell to Text ("@") so Excel does not auto-convert numeric-looking values
# like "589.98" into real numbers, then restore the cell to the default "Normal"
# style so its saved <c> element has no explicit style index, just like the
# original file.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" '69.221.64'
Set-TextValue "E2" '  +2.17%  '

# Row 3
Set-TextValue "D3" '3.314.58'
Set-TextValue "E3" '  +1.94%  '

# Row 4
Set-TextValue "E4" '  -0.01%  '

# Row 5
Set-TextValue "D5" '589.98'
Set-TextValue "E5" '  +1.94%  '

# Row 6
Set-TextValue "D6" '186.66'
Set-TextValue "E6" '  +2.18%  '

# Row 7
Set-TextValue "E7" '  -0.02%  '

# Row 8
Set-TextValue "D8" '0.606'
Set-TextValue "E8" '  +2.11%  '

# Row 9
Set-TextValue "D9" '0.138'
Set-TextValue "E9" '  +5.72%  '

# Row 10
Set-TextValue "D10" '6.72'
Set-TextValue "E10" '  -1.40%  '

# Row 11
Set-TextValue "D11" '0.424'
Set-TextValue "E11" '  +2.54%  '

# Row 12
Set-TextValue "D12" '3.896.05'
Set-TextValue "E12" '  +2.33%  '

# Row 13
Set-TextValue "E13" '  +0.07%  '

# Row 14
Set-TextValue "D14" '29.40'
Set-TextValue "E14" '  +4.71%  '

# Row 15
Set-TextValue "D15" '69.260.44'
Set-TextValue "E15" '  +2.26%  '

# Row 16
Set-TextValue "D16" '0.0000175'
Set-TextValue "E16" '  +3.92%  '

# Row 17
Set-TextValue "D17" '3.314.27'
Set-TextValue "E17" '  +1.90%  '

# Row 18
Set-TextValue "D18" '5.92'
Set-TextValue "E18" '  +1.41%  '

# Row 19
Set-TextValue "D19" '13.83'
Set-TextValue "E19" '  +3.00%  '

# Row 20
Set-TextValue "D20" '393.31'
Set-TextValue "E20" '  +4.98%  '

# Row 21
Set-TextValue "D21" '7.83'
Set-TextValue "E21" '  +3.01%  '

# Row 22
Set-TextValue "D22" '72.09'
Set-TextValue "E22" '  +1.43%  '

# Row 23
Set-TextValue "E23" '  -0.13%  '

# Row 24
Set-TextValue "E24" '  +3.49%  '

# Row 25
Set-TextValue "D25" '0.521'
Set-TextValue "E25" '  +2.09%  '

# Row 26
Set-TextValue "D26" '9.87'
Set-TextValue "E26" '  +2.44%  '

# Row 27
Set-TextValue "D27" '0.189'
Set-TextValue "E27" '  +4.45%  '

# Row 29
Set-TextValue "D29" '5.91'
Set-TextValue "E29" '  +4.18%  '

# Row 30
Set-TextValue "D30" '2.01'
Set-TextValue "E30" '  +1.76%  '

# Row 31
Set-TextValue "B31" 'Fetch.AI'
Set-TextValue "C31" 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue "D31" '1.33'
Set-TextValue "E31" '  +4.38%  '

# Row 32
Set-TextValue "B32" 'EthereumClassic'
Set-TextValue "C32" 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue "D32" '23.19'
Set-TextValue "E32" '  +2.48%  '

# Row 33
Set-TextValue "D33" '7.24'
Set-TextValue "E33" '  +5.11%  '

# Row 34
Set-TextValue "E34" '  +0.01%  '

# Row 35
Set-TextValue "E35" '  +4.44%  '

# Row 36
Set-TextValue "D36" '163.73'
Set-TextValue "E36" '  -0.35%  '

# Row 37
Set-TextValue "D37" '1.93'
Set-TextValue "E37" '  +3.93%  '

# Row 38
Set-TextValue "D38" '0.843'
Set-TextValue "E38" '  -2.42%  '

# Row 39
Set-TextValue "D39" '26.97'
Set-TextValue "E39" '  +0.59%  '

# Row 40
Set-TextValue "B40" 'Filecoin'
Set-TextValue "C40" 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue "D40" '4.64'
Set-TextValue "E40" '  +4.97%  '

# Row 41
Set-TextValue "B41" 'dogwifhat'
Set-TextValue "C41" 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue "D41" '2.66'
Set-TextValue "E41" '  +2.02%  '

# Row 42
Set-TextValue "D42" '6.69'
Set-TextValue "E42" '  -1.58%  '

# Row 43
Set-TextValue "B43" 'Hedera'
Set-TextValue "C43" 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue "D43" '0.0701'
Set-TextValue "E43" '  +3.55%  '

# Row 44
Set-TextValue "B44" 'OKB'
Set-TextValue "C44" 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue "D44" '41.90'
Set-TextValue "E44" '  +3.20%  '

# Row 45
Set-TextValue "B45" 'InjectiveProtocol'
Set-TextValue "C45" 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue "D45" '25.87'
Set-TextValue "E45" '  +0.38%  '

# Row 46
Set-TextValue "D46" '2.672.00'
Set-TextValue "E46" '  -1.36%  '

# Row 47
Set-TextValue "D47" '343.58'
Set-TextValue "E47" '  -5.38%  '

# Row 48
Set-TextValue "E48" '  +3.11%  '

# Row 49
Set-TextValue "D49" '32.70'
Set-TextValue "E49" '  +5.81%  '

# Row 50
Set-TextValue "E50" '  +0.87%  '

# Row 51
Set-TextValue "D51" '6.35'
Set-TextValue "E51" '  +3.71%  '
